# Apply the updated crypto price/volume snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.252.79'
$ws.Range("D3").Value = '2.058.59'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''232.90'
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").Value = '''0.621'
$ws.Range("E6").Value = '  +3.25%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''57.22'
$ws.Range("E8").Value = '  +4.02%  '
$ws.Range("D9").Value = '''0.383'
$ws.Range("E9").Value = '  +3.39%  '
$ws.Range("D10").Value = '''57.76'
$ws.Range("E10").Value = '  +1.46%  '
$ws.Range("D11").Value = '''0.0758'
$ws.Range("E11").Value = '  +0.86%  '
$ws.Range("E12").Value = '  +1.11%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '''14.85'
$ws.Range("E13").Value = '  +3.84%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.359.64'
$ws.Range("E14").Value = '  +1.14%  '
$ws.Range("D15").Value = '''20.91'
$ws.Range("E15").Value = '  +4.43%  '
$ws.Range("D16").Value = '''0.784'
$ws.Range("E16").Value = '  +3.32%  '
$ws.Range("D17").Value = '''5.15'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '2.056.59'
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("D19").Value = '37.192.63'
$ws.Range("E19").Value = '  +1.40%  '
$ws.Range("D20").Value = '''6.35'
$ws.Range("E20").Value = '  +9.02%  '
$ws.Range("D21").Value = '''69.27'
$ws.Range("E21").Value = '  +2.64%  '
$ws.Range("D22").Value = '0.0₃0809'
$ws.Range("E22").Value = '  +1.65%  '
$ws.Range("D23").Value = '''225.39'
$ws.Range("E23").Value = '  +2.23%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '''2.39'
$ws.Range("E26").Value = '  +0.80%  '
$ws.Range("D27").Value = '''165.99'
$ws.Range("E27").Value = '  +2.01%  '
$ws.Range("D28").Value = '''1.46'
$ws.Range("E28").Value = '  +8.48%  '
$ws.Range("D29").Value = '''8.79'
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '''0.127'
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '''19.01'
$ws.Range("E31").Value = '  +0.47%  '
$ws.Range("D32").Value = '''0.117'
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("D33").Value = '''4.46'
$ws.Range("E33").Value = '  +2.28%  '
$ws.Range("D34").Value = '''4.62'
$ws.Range("E34").Value = '  +8.35%  '
$ws.Range("D35").Value = '''0.0618'
$ws.Range("E35").Value = '  +2.24%  '
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("E38").Value = '  +2.22%  '
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("D40").Value = '''5.71'
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("D41").Value = '''2.95'
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.475.64'
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("B43").Value = 'FTXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D43").Value = '''4.38'
$ws.Range("E43").Value = '  -2.94%  '
$ws.Range("D44").Value = '''96.69'
$ws.Range("E44").Value = '  +3.39%  '
$ws.Range("E45").Value = '  +5.23%  '
$ws.Range("D46").Value = '''0.0932'
$ws.Range("E46").Value = '  -0.81%  '
$ws.Range("D47").Value = '''0.0211'
$ws.Range("E47").Value = '  +3.42%  '
$ws.Range("E48").Value = '  +1.95%  '
$ws.Range("D49").Value = '''15.10'
$ws.Range("E49").Value = '  -2.92%  '
$ws.Range("D50").Value = '''7.16'
$ws.Range("E50").Value = '  +3.02%  '
$ws.Range("E51").Value = '  +1.79%  '
